$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '41.774.37'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '2.468.71'
$ws.Range('E3').Value = '  -0.62%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.15%  '
Set-TextValue $ws.Range('D5') '316.30'
$ws.Range('E5').Value = '  +1.22%  '
Set-TextValue $ws.Range('D6') '93.02'
$ws.Range('E6').Value = '  +0.02%  '
Set-TextValue $ws.Range('D7') '0.549'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E9').Value = '  +3.72%  '
Set-TextValue $ws.Range('D10') '32.87'
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  +8.83%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '2.849.25'
$ws.Range('E13').Value = '  -0.53%  '
Set-TextValue $ws.Range('D14') '6.90'
$ws.Range('E14').Value = '  +0.63%  '
Set-TextValue $ws.Range('D15') '15.77'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '2.469.41'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').Value = '41.754.05'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D21') '11.55'
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D22') '71.12'
$ws.Range('E22').Value = '  +0.85%  '
Set-TextValue $ws.Range('D23') '239.35'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('E26').Value = '  -0.07%  '
Set-TextValue $ws.Range('D27') '24.87'
$ws.Range('E27').Value = '  +0.37%  '
Set-TextValue $ws.Range('D28') '2.26'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('E30').Value = '  -0.97%  '
Set-TextValue $ws.Range('D31') '155.88'
$ws.Range('E31').Value = '  +0.81%  '
Set-TextValue $ws.Range('D32') '5.54'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D33') '2.58'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D34') '0.0766'
$ws.Range('E34').Value = '  +1.65%  '
$ws.Range('B35').Value = 'ApeXProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D35') '2.49'
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D36') '17.59'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('E38').Value = '  +1.44%  '
Set-TextValue $ws.Range('D39') '0.104'
$ws.Range('E39').Value = '  -1.31%  '
Set-TextValue $ws.Range('D41') '4.00'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '1.973.09'
$ws.Range('E43').Value = '  +0.75%  '
Set-TextValue $ws.Range('D44') '18.96'
$ws.Range('E44').Value = '  -4.12%  '
$ws.Range('E45').Value = '  -0.38%  '
Set-TextValue $ws.Range('D46') '2.95'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('D48').Value = '2.701.78'
$ws.Range('E48').Value = '  -0.89%  '
Set-TextValue $ws.Range('D49') '96.85'
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('E50').Value = '  +0.26%  '
Set-TextValue $ws.Range('D51') '72.86'
$ws.Range('E51').Value = '  -0.43%  '
